# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and
# WVR sheets with freshly pulled prices, and recompute the profit figures
# that depend on them. A few rows gain or lose a stray HQ/NQ profit cell
# where the sourced item no longer has (or newly has) an HQ/NQ variant.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 36.2
$ws.Range("I8").Value = 36.2
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 108.6
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 30.39999999999999

$ws.Range("H107").Value = 486.94116
$ws.Range("I107").Value = 321.18182
$ws.Range("J107").Value = 790.8333
$ws.Range("K107").Value = 321.18182
$ws.Range("L107").Value = 790.8333
$ws.Range("M107").Value = 1598.81818
$ws.Range("N107").Value = -4630.8333

$ws.Range("H116").Value = 2327.1667
$ws.Range("I116").Value = 2115
$ws.Range("J116").Value = 2751.5
$ws.Range("K116").Value = 2115
$ws.Range("L116").Value = 2751.5
$ws.Range("M116").Value = 1327
$ws.Range("N116").Value = -9635.5

$ws.Range("H129").Value = 1070.3276
$ws.Range("I129").Value = 2250.3333
$ws.Range("J129").Value = 1005.9636
$ws.Range("K129").Value = 6750.999899999999
$ws.Range("L129").Value = 3017.8908
$ws.Range("M129").Value = -1750.999899999999
$ws.Range("N129").Value = -13017.8908

$ws.Range("H132").Value = 10876260
$ws.Range("I132").Value = 14291936
$ws.Range("J132").Value = 8199
$ws.Range("K132").Value = 42875808
$ws.Range("L132").Value = 24597
$ws.Range("M132").Value = -42873278
$ws.Range("N132").Value = -29657

$ws.Range("H138").Value = 1197.0513
$ws.Range("I138").Value = 699.5161000000001
$ws.Range("J138").Value = 3125
$ws.Range("K138").Value = 2098.5483
$ws.Range("L138").Value = 9375
$ws.Range("M138").Value = 3041.4517
$ws.Range("N138").Value = -19655

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -827

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("L17").ClearContents()

$ws.Range("H74").Value = 1699.75
$ws.Range("I74").Value = 1699.75
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1699.75
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = -825.75
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 1699.75
$ws.Range("I77").Value = 1699.75
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8498.75
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = -4130.75
$ws.Range("M77").ClearContents()

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -3550
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 1509669.6
$ws.Range("I132").Value = 992.85187
$ws.Range("J132").Value = 4904192.5
$ws.Range("K132").Value = 2978.55561
$ws.Range("L132").Value = 14712577.5
$ws.Range("M132").Value = -448.5556099999999
$ws.Range("N132").Value = -14717637.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2472318.8
$ws.Range("I134").Value = 716.45715
$ws.Range("J134").Value = 11122927
$ws.Range("K134").Value = 2149.37145
$ws.Range("L134").Value = 33368781
$ws.Range("M134").Value = 385.6285500000004
$ws.Range("N134").Value = -33373851

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 136.42857
$ws.Range("I19").Value = 136.42857
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 136.42857
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = 33.57142999999999
$ws.Range("M19").ClearContents()

$ws.Range("H24").Value = 136.42857
$ws.Range("I24").Value = 136.42857
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 136.42857
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = 33.57142999999999
$ws.Range("M24").ClearContents()

$ws.Range("H62").Value = 4175.4443
$ws.Range("I62").Value = 2899.75
$ws.Range("J62").Value = 5196
$ws.Range("K62").Value = 2899.75
$ws.Range("L62").Value = 5196
$ws.Range("M62").Value = -2275.75
$ws.Range("N62").Value = -6444

$ws.Range("H65").Value = 4175.4443
$ws.Range("I65").Value = 2899.75
$ws.Range("J65").Value = 5196
$ws.Range("K65").Value = 14498.75
$ws.Range("L65").Value = 25980
$ws.Range("M65").Value = -11378.75
$ws.Range("N65").Value = -32220

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("N116").Value = 0
$ws.Range("L116").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 123.111115
$ws.Range("I10").Value = 123.111115
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 369.333345
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -230.333345

$ws.Range("H22").Value = 9105489
$ws.Range("I22").Value = 100000000
$ws.Range("J22").Value = 16037.4
$ws.Range("K22").Value = 300000000
$ws.Range("L22").Value = 48112.2
$ws.Range("M22").Value = -299999831
$ws.Range("N22").Value = -48450.2

$ws.Range("H27").Value = 9105489
$ws.Range("I27").Value = 100000000
$ws.Range("J27").Value = 16037.4
$ws.Range("K27").Value = 300000000
$ws.Range("L27").Value = 48112.2
$ws.Range("M27").Value = -299999898
$ws.Range("N27").Value = -48316.2

$ws.Range("H40").Value = 403.55554
$ws.Range("I40").Value = 166.4
$ws.Range("J40").Value = 700
$ws.Range("K40").Value = 665.6
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -596.6
$ws.Range("N40").Value = -2938

$ws.Range("H129").Value = 1351.3
$ws.Range("I129").Value = 926
$ws.Range("J129").Value = 1493.0667
$ws.Range("K129").Value = 2778
$ws.Range("L129").Value = 4479.2001
$ws.Range("M129").Value = 2222
$ws.Range("N129").Value = -14479.2001

$ws.Range("H131").Value = 838.09
$ws.Range("I131").Value = 423.75
$ws.Range("J131").Value = 874.11957
$ws.Range("K131").Value = 1271.25
$ws.Range("L131").Value = 2622.35871
$ws.Range("M131").Value = 3768.75
$ws.Range("N131").Value = -12702.35871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 3998
$ws.Range("I28").Value = 3998
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3998
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -3806
$ws.Range("N28").ClearContents()

$ws.Range("H97").Value = 665.2162
$ws.Range("I97").Value = 548.4706
$ws.Range("J97").Value = 764.45
$ws.Range("K97").Value = 548.4706
$ws.Range("L97").Value = 764.45
$ws.Range("M97").Value = -52.47059999999999
$ws.Range("N97").Value = -1756.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 8646.666999999999
$ws.Range("I17").Value = 2970
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 2970
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = -2800
$ws.Range("N17").Value = -20340

$ws.Range("H24").Value = 30000000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 30000000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 30000000
$ws.Range("N24").Value = -30000686

$ws.Range("H30").Value = 1436.4
$ws.Range("I30").Value = 545.5
$ws.Range("J30").Value = 5000
$ws.Range("K30").Value = 545.5
$ws.Range("L30").Value = 5000
$ws.Range("M30").Value = -437.5
$ws.Range("N30").Value = -5216

$ws.Range("H55").Value = 62506520
$ws.Range("I55").Value = 20148.2
$ws.Range("J55").Value = 90909416
$ws.Range("K55").Value = 20148.2
$ws.Range("L55").Value = 90909416
$ws.Range("M55").Value = -19975.2
$ws.Range("N55").Value = -90909762

$ws.Range("H136").Value = 33615356
$ws.Range("I136").Value = 4763638.5
$ws.Range("J136").Value = 250003250
$ws.Range("K136").Value = 14290915.5
$ws.Range("L136").Value = 750009750
$ws.Range("M136").Value = -14288365.5
$ws.Range("N136").Value = -750014850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 3000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3470

$ws.Range("H35").Value = 3000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 3000
$ws.Range("N35").Value = -3580

$ws.Range("H107").Value = 610.1667
$ws.Range("I107").Value = 677.375
$ws.Range("J107").Value = 475.75
$ws.Range("K107").Value = 2032.125
$ws.Range("L107").Value = 1427.25
$ws.Range("M107").Value = -112.125
$ws.Range("N107").Value = -5267.25
